# Build site at 2023-04-12 14:53:07 UTC
# Fix misaligned label/value rows on the LOQ4005 sheet:
#  - "Objetivos:" (row 10) gets its own real objectives text instead of the
#    teacher-name text that had been pasted there by mistake.
#  - A new row is inserted under "Docentes responsáveis:" to hold the
#    teacher-name value, and every following value shifts down into the row
#    that matches its real label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at 13 (pushes old rows 13-23 down to 14-24) ---
$ws.Rows(13).Insert()

# The inserted row inherited column A's bold style from the row above
# ("Docentes responsáveis:"); this row should have no label in column A,
# only the teacher name in B/C, so copy the B/C formatting from row 10 and
# drop the stray A13 formatting.
$ws.Range("B10:C10").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Cells.Item(13, 1).Clear()

# --- Fix the column definitions: column 1 should only cover column A ---
# (the original file incorrectly grouped columns A and B into one <col>
# definition with A's width/style; nudging column B's width splits that
# group so column A gets its own, correctly scoped definition)
$ws.Columns(2).ColumnWidth = 60.7109375

# --- Correct the cell values so each label lines up with its real value ---

$ws.Range("B10:C10").Value = "1) Gerais: - Levar os estudantes a compreenderem os mecanismos de obtenção e análise de daos de variáveis de processo ,identificando as causas especiais de variação ( instabilidade), e causas comuns, de natureza aleatória. 2) Específicos: - Ao final do curso os educandos devem:? Saber identificar causas especiais de variação? Saber analisar os resultados propondo a condições que levem o processo a uma menor variabilidadeSaber determinar a capacidade do processo e utilizar as cartas de controle"

$ws.Range("B13:C13").Value = "5840535 - Messias Borges Silva"

$ws.Range("B14:C14").Value = "IntroduçãoCartas de ControleAnálise da Cacidade de ProcessosCartas EspeciasCasos Práticos"

$ws.Range("B16:C16").Value = "I - DESCRITIVO:INTRODUÇÃO- A importância do CEP- Potencialidades- Natureza da variação - Causas especiais e causas comunsCARTAS DE CONTROLE- Cartas X, R- Carta X- Carta P- Carta nP- Carta C- Carta UCARTAS DE CONTROLE ESPECIAIS- Amplitude móvel- Soma acumulada (CUSUM)ANÁLISE DE CAPACIDADE DOS PROCESSOS- Indice Co- Indice CpK- Indice PPKCASOS PRÁTICOS- Utilização de situações reais vivenciados em ambiente indus-trial."

$ws.Range("B19:C19").Value = "duas provas escritas"

$ws.Range("B20:C20").Value = "serão avaliados os conteúdos discutidos em sala e constantes da ementa do curso.A média da disciplina será a média aritmética das duas provas."

$ws.Range("B21:C21").Value = "uma prova escrita com conteúdo de todo o semestre"

$ws.Range("B22:C22").Value = "1 - Ferramentas Estatísticas Básicas p/ o Gerenciamento de Processos. Maria Cristina C. Werkena. Edit. FCO, 19962) Controle Estatístico da Qualidade, 4ª edição. Douglas C. Mont gomery, 2006. Edit. LTC3- Statisticial Quality Control, 5ª edição. Fugeno L. Grant and Richard S. Leavenworth MC Graw Hill, 1987"
